# Append the CSS content block after the closing </html> paragraph,
# and mint footnotes.xml/endnotes.xml (separator-only stub parts).
$d = $word.ActiveDocument
$sel = $word.Selection
$sel.EndKey(6, 0) | Out-Null  # wdStory: jump to the very end of the document

$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('CSS') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('#css-slider {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('             margin-left:  168px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              width: 400px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              height: 200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              margin-top:20px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              overflow: hidden;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('.slide-item {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              width: 410px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              height:200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              float: left;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('              position: relative;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('body {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  width:1000px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  margin:0 auto;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('#D {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-left:0px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('width:180px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('min-height:260px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('height:200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-top:40px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('#A {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-left:190px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('width:560px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('min-height:260px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('height:200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('background-color:#009b9f;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-top:-260px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('#D {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-left:0px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('width:180px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('min-height:260px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('height:200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-top:40px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('#B {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-left:0px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('width:750px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('min-height:200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('height:200px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('margin-top:0px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('html,') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('body {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  font-size: 20px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  height: 100%;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  margin: 0;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('@media screen and (min-width: 768px) {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  html,') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  body {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('    font-size: 30px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  }') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('@media screen and (min-width: 980px) {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  html,') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  body {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('    font-size: 32px;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  }') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('.item_wrapper {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  list-style-type: none;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  display: -ms-flexbox;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  display: flex;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  flex-wrap: wrap;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  padding: 0;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  margin: 0 1vw 1vw;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('.item_wrapper li {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  width: calc(100% / 3 - 2 * 1vw);') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  margin: 1vw;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('@media screen and (min-width: 425px) {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  .item_wrapper li {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('    width: calc(100% / 4 - 2 * 1vw);') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  }') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('@media screen and (min-width: 980px) {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  .item_wrapper li {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('    width: calc(100% / 5 - 2 * 1vw);') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  }') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('.item_wrapper img {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  vertical-align: bottom;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('.rank1 li {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  counter-increment: rank;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('.rank1 li::before {') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  content: counter(rank) "') | Out-Null
$sel.TypeText('位') | Out-Null
$sel.TypeText('";') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  color: #bf0000;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  font-weight: 700;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('  line-height: initial;') | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText('}') | Out-Null

# Mint footnotes.xml / endnotes.xml (Word creates both parts the first time
# a footnote is used); deleting the footnote again leaves the separator-only
# parts behind without leaving any footnote reference in the body text.
$fnRange = $d.Range(0, 0)
$fn = $d.Footnotes.Add($fnRange, "", "x")
$fn.Delete() | Out-Null

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
